$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2..38) are being rearranged: the full record (all 18 columns)
# that currently lives in a given row is relocated to another row. Capture the
# whole block first, then write it back out in the new order so that no values
# are lost while rows are being permuted.

$range = $ws.Range("A2:R38")
$data = $range.Value2

# Map: new row position (1 = row2 ... 37 = row38) -> source row position in $data
$order = @(37,9,32,23,11,27,16,21,29,19,1,10,33,35,30,15,2,31,8,14,7,18,6,4,25,22,13,36,20,3,12,26,28,5,24,17,34)

$rows = $range.Rows.Count
$cols = $range.Columns.Count

$new = New-Object 'object[,]' $rows,$cols
for ($i = 1; $i -le $rows; $i++) {
    $srcRow = $order[$i-1]
    for ($j = 1; $j -le $cols; $j++) {
        $new[$i-1, $j-1] = $data[$srcRow, $j]
    }
}

$range.Value2 = $new
